$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.928.04"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "3.103.27"
$ws.Range("E3").Value = "  +5.21%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "580.99"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "173.02"
$ws.Range("E6").Value = "  +6.82%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.099.08"
$ws.Range("E8").Value = "  +5.20%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "6.46"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("E12").Value = "  +4.14%  "
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "37.31"
$ws.Range("E14").Value = "  +7.20%  "
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "3.617.26"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "66.897.62"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "3.109.51"
$ws.Range("E19").Value = "  +5.38%  "
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "482.67"
$ws.Range("E21").Value = "  +8.25%  "
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("D24").Value = "84.03"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  +6.49%  "
$ws.Range("D26").Value = "13.12"
$ws.Range("E26").Value = "  +7.13%  "
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -4.63%  "
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "28.73"
$ws.Range("E33").Value = "  +5.92%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("D37").Value = "0.995"
$ws.Range("E37").Value = "  +2.23%  "
$ws.Range("D38").Value = "48.18"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("D41").Value = "50.15"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("D46").Value = "2.832.61"
$ws.Range("E46").Value = "  +5.61%  "
$ws.Range("D47").Value = "384.06"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "134.87"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "24.87"
$ws.Range("E50").Value = "  +4.33%  "
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  +3.15%  "
